# Refresh the cryptos Price/Volume(1h) columns (and swap the HuobiToken /
# TrustWalletToken rows) per the latest GitHub Actions data pull.
# Note: a leading '' (single-quote) is used on Price values that Excel would
# otherwise auto-parse as a number, so they stay plain text like the rest of
# column D (matching values such as "27.151.56" that already contain
# multiple dots and can only ever be text).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.151.56'
$ws.Range("E2").Value = '  -1.87%  '

$ws.Range("D3").Value = '1.563.17'
$ws.Range("E3").Value = '  -1.58%  '

$ws.Range("E4").Value = '  -0.08%  '

$ws.Range("D5").Value = '''206.64'
$ws.Range("E5").Value = '  -0.19%  '

$ws.Range("E6").Value = '  -1.71%  '

$ws.Range("E7").Value = '  -0.05%  '

$ws.Range("E8").Value = '  +0.13%  '

$ws.Range("E9").Value = '  -2.08%  '

$ws.Range("E10").Value = '  +0.25%  '

$ws.Range("D11").Value = '''0.0861'
$ws.Range("E11").Value = '  -0.63%  '

$ws.Range("D12").Value = '1.784.62'
$ws.Range("E12").Value = '  -1.63%  '

$ws.Range("D13").Value = '1.555.48'
$ws.Range("E13").Value = '  -2.07%  '

$ws.Range("E14").Value = '  -2.04%  '

$ws.Range("D15").Value = '''0.517'
$ws.Range("E15").Value = '  -2.42%  '

$ws.Range("D16").Value = '''62.87'
$ws.Range("E16").Value = '  -0.87%  '

$ws.Range("D17").Value = '27.152.62'
$ws.Range("E17").Value = '  -1.81%  '

$ws.Range("D18").Value = '''213.08'
$ws.Range("E18").Value = '  -2.77%  '

$ws.Range("E19").Value = '  -1.20%  '

$ws.Range("E20").Value = '  -1.39%  '

$ws.Range("E21").Value = '  -0.03%  '

$ws.Range("D22").Value = '''4.12'
$ws.Range("E22").Value = '  -0.35%  '

$ws.Range("D23").Value = '''9.40'
$ws.Range("E23").Value = '  -1.90%  '

$ws.Range("E24").Value = '  +0.24%  '

$ws.Range("D25").Value = '''152.15'
$ws.Range("E25").Value = '  -0.95%  '

$ws.Range("D26").Value = '''6.58'
$ws.Range("E26").Value = '  -3.75%  '

$ws.Range("E27").Value = '  -1.48%  '

$ws.Range("E29").Value = '  -1.39%  '

$ws.Range("E30").Value = '  -0.91%  '

$ws.Range("D31").Value = '''0.0463'
$ws.Range("E31").Value = '  -0.87%  '

$ws.Range("E32").Value = '  -1.76%  '

$ws.Range("D33").Value = '1.382.26'
$ws.Range("E33").Value = '  +0.80%  '

$ws.Range("E34").Value = '  +0.57%  '

$ws.Range("E35").Value = '  +0.23%  '

$ws.Range("B36").Value = 'TrustWalletToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D36").Value = '''0.945'
$ws.Range("E36").Value = '  -3.52%  '

$ws.Range("B37").Value = 'HuobiToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D37").Value = '''2.28'
$ws.Range("E37").Value = '  -1.01%  '

$ws.Range("E38").Value = '  -1.13%  '

$ws.Range("D39").Value = '''0.815'
$ws.Range("E39").Value = '  -1.13%  '

$ws.Range("E40").Value = '  -3.31%  '

$ws.Range("E42").Value = '  +1.75%  '

$ws.Range("D43").Value = '''1.79'
$ws.Range("E43").Value = '  +4.04%  '

$ws.Range("E44").Value = '  +0.01%  '

$ws.Range("D45").Value = '''63.36'
$ws.Range("E45").Value = '  -1.12%  '

$ws.Range("E46").Value = '  +0.28%  '

$ws.Range("D47").Value = '1.697.85'
$ws.Range("E47").Value = '  -1.60%  '

$ws.Range("D48").Value = '''85.64'
$ws.Range("E48").Value = '  -2.10%  '

$ws.Range("D49").Value = '0.0₇0993'
$ws.Range("E49").Value = '  -0.99%  '

$ws.Range("D50").Value = '''0.0492'
$ws.Range("E50").Value = '  -0.50%  '

$ws.Range("E51").Value = '  +0.12%  '
